{"js": "const body = context.document.body;\nconst paras = body.paragraphs;\nparas.load(\"items\");\nawait context.sync();\nconst heading = paras.items[22];\nconst newPara = heading.insertParagraph(undefined, \"Before\");\nawait context.sync();\nnewPara.getRange().clear();\nawait context.sync();\n", "ps1": "$d = $word.ActiveDocument\n$p = $d.Paragraphs.Item(23)\n$r = $p.Range.Duplicate\n$r.Collapse(1)\n$r.InsertParagraphBefore()\n"}
